# The document contains two (near-duplicate) instances of the instruction
# "Create a MovePenguinToStart function for the Penguin 2 sprite." The
# second instance (the one right after the page-break, further down in the
# document) is the one whose wording changed to "MovePenguin" and which
# picked up Word's auto-managed "_GoBack" bookmark (Word moves "_GoBack" to
# mark the most recent edit location). The old "_GoBack" bookmark, which
# used to sit around the final picture near the end of the document, is
# removed as a side effect of the bookmark being moved.

$d = $word.ActiveDocument

# --- Locate the first occurrence of "MovePenguinToStart" ---
$rngFirst = $d.Content
$rngFirst.Find.ClearFormatting()
$rngFirst.Find.Execute("MovePenguinToStart", $true, $false, $false, $false, `
                        $false, $true, 1, $false, "", 0) | Out-Null

# --- Locate the second occurrence, searching after the first one ---
$rngSecond = $d.Range($rngFirst.End, $d.Content.End)
$rngSecond.Find.ClearFormatting()
$rngSecond.Find.Execute("MovePenguinToStart", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "", 0) | Out-Null

# --- Rename it to "MovePenguin" (only this, second, occurrence) ---
$rngSecond.Text = "MovePenguin"

# Collapse the range to the end of the freshly written text, i.e. right
# after "MovePenguin", so the bookmark we add next is zero-length and
# sits immediately after that run (matching where Word leaves "_GoBack"
# after the user's last edit).
$rngSecond.Collapse(0)

# --- Move the "_GoBack" bookmark here ---
# Adding a bookmark with a name that already exists elsewhere in the
# document relocates it (Word keeps bookmark names unique), which removes
# the old "_GoBack" bookmark that used to wrap the final picture.
$d.Bookmarks.Add("_GoBack", $rngSecond) | Out-Null

Write-Output "MovePenguinToStart -> MovePenguin renamed and _GoBack bookmark relocated"
